$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.261.82'
$ws.Range("E2").Value = '  +3.93%  '

$ws.Range("D3").Value = '1.787.61'
$ws.Range("E3").Value = '  +0.20%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '339.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3832'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3447'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.02'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.156'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07418'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.30'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.90%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.9984'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.471'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.50%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.377'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.02%  '

$ws.Range("D16").Value = '1.785.32'
$ws.Range("E16").Value = '  +0.40%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001079'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06682'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.483'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.22%  '

$ws.Range("D23").Value = '28.299.78'
$ws.Range("E23").Value = '  +4.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.355'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.85%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.452'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.37%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.430'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.59'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.27%  '

$ws.Range("D31").Value = '1.983.31'
$ws.Range("E31").Value = '  +0.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.150'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.65%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.955'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.51%  '

$ws.Range("E34").Value = '  +2.72%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.83'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.42%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02432'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.87%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6889'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.364'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.67%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06393'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.43%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2177'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.93%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.249'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.500'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.87%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.305'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.74%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.18'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9995'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6322'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.95%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.877'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.51%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.089'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07492'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.213'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.37%  '
